# Automatic update of files.
# Rows 28 and 29 swap their data (the "Spillkråka" record and the
# "Tretåig hackspett" record exchange positions), and four brand new
# observation rows (30-33) are appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 28 <- previous row 29 content (Tretåig hackspett / hane / födosökande)
# ---------------------------------------------------------------------
$ws.Cells.Item(28, 1).Value = 131090275
$ws.Cells.Item(28, 2).Value = 57884
$ws.Cells.Item(28, 4).Value = "NT"
$ws.Cells.Item(28, 5).Value = 100109
$ws.Cells.Item(28, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(28, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(28, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(28, 12).Value = "hane"
$ws.Cells.Item(28, 13).Value = "födosökande"
$ws.Cells.Item(28, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(28, 17).Value = 584987
$ws.Cells.Item(28, 18).Value = 7060190
$ws.Cells.Item(28, 19).Value = 15
$ws.Cells.Item(28, 20).Value = "Västernorrland"
$ws.Cells.Item(28, 21).Value = "Sollefteå"
$ws.Cells.Item(28, 22).Value = "Ångermanland"
$ws.Cells.Item(28, 23).Value = "Junsele"
$ws.Cells.Item(28, 25).Value = "'2026-02-09"
$ws.Cells.Item(28, 26).ClearContents()
$ws.Cells.Item(28, 27).Value = "'2026-02-09"
$ws.Cells.Item(28, 28).ClearContents()
$ws.Cells.Item(28, 30).Value = $false
$ws.Cells.Item(28, 31).Value = $false
$ws.Cells.Item(28, 33).Value = $false
$ws.Cells.Item(28, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(28, 50).Value = "Daniel Rutschman"

# ---------------------------------------------------------------------
# Row 29 <- previous row 28 content (Spillkråka / färska spår)
# ---------------------------------------------------------------------
$ws.Cells.Item(29, 1).Value = 131089521
$ws.Cells.Item(29, 2).Value = 57881
$ws.Cells.Item(29, 4).Value = "NT"
$ws.Cells.Item(29, 5).Value = 100049
$ws.Cells.Item(29, 6).Value = "Spillkråka"
$ws.Cells.Item(29, 7).Value = "Dryocopus martius"
$ws.Cells.Item(29, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(29, 12).ClearContents()
$ws.Cells.Item(29, 13).Value = "färska spår"
$ws.Cells.Item(29, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(29, 17).Value = 584995
$ws.Cells.Item(29, 18).Value = 7060537
$ws.Cells.Item(29, 19).Value = 15
$ws.Cells.Item(29, 20).Value = "Västernorrland"
$ws.Cells.Item(29, 21).Value = "Sollefteå"
$ws.Cells.Item(29, 22).Value = "Ångermanland"
$ws.Cells.Item(29, 23).Value = "Junsele"
$ws.Cells.Item(29, 25).Value = "'2026-02-09"
$ws.Cells.Item(29, 26).Value = "'13:14"
$ws.Cells.Item(29, 27).Value = "'2026-02-09"
$ws.Cells.Item(29, 28).Value = "'13:14"
$ws.Cells.Item(29, 30).Value = $false
$ws.Cells.Item(29, 31).Value = $false
$ws.Cells.Item(29, 33).Value = $false
$ws.Cells.Item(29, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(29, 50).Value = "Daniel Rutschman"

# ---------------------------------------------------------------------
# Row 30 (new): Tretåig hackspett, färska spår, Röån
# ---------------------------------------------------------------------
$ws.Cells.Item(30, 1).Value = 131144497
$ws.Cells.Item(30, 2).Value = 57884
$ws.Cells.Item(30, 4).Value = "NT"
$ws.Cells.Item(30, 5).Value = 100109
$ws.Cells.Item(30, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(30, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(30, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(30, 13).Value = "färska spår"
$ws.Cells.Item(30, 16).Value = "Röån, Ång"
$ws.Cells.Item(30, 17).Value = 584884
$ws.Cells.Item(30, 18).Value = 7060423
$ws.Cells.Item(30, 19).Value = 15
$ws.Cells.Item(30, 20).Value = "Västernorrland"
$ws.Cells.Item(30, 21).Value = "Sollefteå"
$ws.Cells.Item(30, 22).Value = "Ångermanland"
$ws.Cells.Item(30, 23).Value = "Junsele"
$ws.Cells.Item(30, 25).Value = "'2026-02-09"
$ws.Cells.Item(30, 27).Value = "'2026-02-09"
$ws.Cells.Item(30, 29).Value = "Färska ringhack"
$ws.Cells.Item(30, 30).Value = $false
$ws.Cells.Item(30, 31).Value = $false
$ws.Cells.Item(30, 33).Value = $false
$ws.Cells.Item(30, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(30, 50).Value = "Daniel Rutschman"

# ---------------------------------------------------------------------
# Row 31 (new): Tretåig hackspett, hane-less, födosökande, Antal=1
# ---------------------------------------------------------------------
$ws.Cells.Item(31, 1).Value = 131144494
$ws.Cells.Item(31, 2).Value = 57884
$ws.Cells.Item(31, 4).Value = "NT"
$ws.Cells.Item(31, 5).Value = 100109
$ws.Cells.Item(31, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(31, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(31, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(31, 9).Value = "'1"
$ws.Cells.Item(31, 13).Value = "födosökande"
$ws.Cells.Item(31, 16).Value = "Röån, Ång"
$ws.Cells.Item(31, 17).Value = 585030
$ws.Cells.Item(31, 18).Value = 7060258
$ws.Cells.Item(31, 19).Value = 15
$ws.Cells.Item(31, 20).Value = "Västernorrland"
$ws.Cells.Item(31, 21).Value = "Sollefteå"
$ws.Cells.Item(31, 22).Value = "Ångermanland"
$ws.Cells.Item(31, 23).Value = "Junsele"
$ws.Cells.Item(31, 25).Value = "'2026-02-09"
$ws.Cells.Item(31, 27).Value = "'2026-02-09"
$ws.Cells.Item(31, 30).Value = $false
$ws.Cells.Item(31, 31).Value = $false
$ws.Cells.Item(31, 33).Value = $false
$ws.Cells.Item(31, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(31, 50).Value = "Daniel Rutschman"

# ---------------------------------------------------------------------
# Row 32 (new): Tretåig hackspett, färska spår, Röån (2nd spot)
# ---------------------------------------------------------------------
$ws.Cells.Item(32, 1).Value = 131144496
$ws.Cells.Item(32, 2).Value = 57884
$ws.Cells.Item(32, 4).Value = "NT"
$ws.Cells.Item(32, 5).Value = 100109
$ws.Cells.Item(32, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(32, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(32, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(32, 13).Value = "färska spår"
$ws.Cells.Item(32, 16).Value = "Röån, Ång"
$ws.Cells.Item(32, 17).Value = 584875
$ws.Cells.Item(32, 18).Value = 7060422
$ws.Cells.Item(32, 19).Value = 15
$ws.Cells.Item(32, 20).Value = "Västernorrland"
$ws.Cells.Item(32, 21).Value = "Sollefteå"
$ws.Cells.Item(32, 22).Value = "Ångermanland"
$ws.Cells.Item(32, 23).Value = "Junsele"
$ws.Cells.Item(32, 25).Value = "'2026-02-09"
$ws.Cells.Item(32, 27).Value = "'2026-02-09"
$ws.Cells.Item(32, 29).Value = "Färska ringhack, tall"
$ws.Cells.Item(32, 30).Value = $false
$ws.Cells.Item(32, 31).Value = $false
$ws.Cells.Item(32, 33).Value = $false
$ws.Cells.Item(32, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(32, 50).Value = "Daniel Rutschman"

# ---------------------------------------------------------------------
# Row 33 (new): Talltita, lockläte/övriga läten, Röån
# ---------------------------------------------------------------------
$ws.Cells.Item(33, 1).Value = 131144498
$ws.Cells.Item(33, 2).Value = 58043
$ws.Cells.Item(33, 4).Value = "NT"
$ws.Cells.Item(33, 5).Value = 103021
$ws.Cells.Item(33, 6).Value = "Talltita"
$ws.Cells.Item(33, 7).Value = "Poecile montanus"
$ws.Cells.Item(33, 8).Value = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(33, 13).Value = "lockläte, övriga läten"
$ws.Cells.Item(33, 16).Value = "Röån, Ång"
$ws.Cells.Item(33, 17).Value = 584857
$ws.Cells.Item(33, 18).Value = 7060494
$ws.Cells.Item(33, 19).Value = 15
$ws.Cells.Item(33, 20).Value = "Västernorrland"
$ws.Cells.Item(33, 21).Value = "Sollefteå"
$ws.Cells.Item(33, 22).Value = "Ångermanland"
$ws.Cells.Item(33, 23).Value = "Junsele"
$ws.Cells.Item(33, 25).Value = "'2026-02-09"
$ws.Cells.Item(33, 27).Value = "'2026-02-09"
$ws.Cells.Item(33, 30).Value = $false
$ws.Cells.Item(33, 31).Value = $false
$ws.Cells.Item(33, 33).Value = $false
$ws.Cells.Item(33, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(33, 50).Value = "Daniel Rutschman"
